$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '44.030.20'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '2.367.11'
$ws.Range("E3").Value = '  +4.93%  '
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue "D5" '235.52'
$ws.Range("E5").Value = '  +1.88%  '
$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue "D6" '0.660'
$ws.Range("E6").Value = '  +3.14%  '
Set-TextValue "D7" '73.70'
$ws.Range("E7").Value = '  +14.66%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-TextValue "D9" '0.531'
$ws.Range("E9").Value = '  +21.14%  '
$ws.Range("E10").Value = '  +3.69%  '
Set-TextValue "D11" '28.45'
$ws.Range("E11").Value = '  +7.77%  '
$ws.Range("D12").Value = '2.715.07'
$ws.Range("E12").Value = '  +4.81%  '
$ws.Range("E13").Value = '  +2.35%  '
Set-TextValue "D14" '16.85'
$ws.Range("E14").Value = '  +12.76%  '
Set-TextValue "D16" '0.886'
$ws.Range("E16").Value = '  +7.78%  '
$ws.Range("D17").Value = '2.356.23'
$ws.Range("E17").Value = '  +4.02%  '
$ws.Range("D18").Value = '43.858.92'
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("E19").Value = '  +4.66%  '
Set-TextValue "D20" '76.06'
$ws.Range("E20").Value = '  +4.48%  '
Set-TextValue "D21" '6.34'
$ws.Range("E21").Value = '  +4.21%  '
Set-TextValue "D22" '251.65'
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("E23").Value = '  +0.03%  '
Set-TextValue "D24" '3.77'
$ws.Range("E24").Value = '  -1.98%  '
Set-TextValue "D25" '2.49'
$ws.Range("E25").Value = '  +1.95%  '
Set-TextValue "D26" '10.27'
$ws.Range("E26").Value = '  +6.05%  '
$ws.Range("E27").Value = '  -1.40%  '
Set-TextValue "D28" '22.55'
$ws.Range("E28").Value = '  +4.31%  '
Set-TextValue "D29" '173.02'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("E30").Value = '  +9.10%  '
$ws.Range("E31").Value = '  +1.92%  '
$ws.Range("E32").Value = '  +4.82%  '
$ws.Range("E33").Value = '  +4.46%  '
Set-TextValue "D34" '0.0709'
$ws.Range("E34").Value = '  +4.73%  '
Set-TextValue "D35" '5.13'
$ws.Range("E35").Value = '  +4.85%  '
Set-TextValue "D36" '3.77'
$ws.Range("E36").Value = '  +4.56%  '
Set-TextValue "D37" '2.45'
$ws.Range("E37").Value = '  +8.45%  '
Set-TextValue "D38" '6.46'
$ws.Range("E38").Value = '  +0.93%  '
$ws.Range("E39").Value = '  +5.85%  '
Set-TextValue "D40" '19.34'
$ws.Range("E40").Value = '  +13.14%  '
$ws.Range("E41").Value = '  +0.08%  '
Set-TextValue "D42" '8.91'
$ws.Range("E42").Value = '  +1.42%  '
Set-TextValue "D43" '1.18'
$ws.Range("E43").Value = '  +10.24%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D44" '1.22'
$ws.Range("E44").Value = '  +3.24%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D45" '98.94'
$ws.Range("E45").Value = '  +2.53%  '
Set-TextValue "D46" '0.0966'
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("E48").Value = '  +13.40%  '
$ws.Range("D49").Value = '1.441.85'
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D50" '2.30'
$ws.Range("E50").Value = '  +1.98%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.588.60'
$ws.Range("E51").Value = '  +4.76%  '
